$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").ClearContents()
$ws.Range("F23").ClearContents()

$ws.Range("H19").Select()
